# Insert a new data row at row 458 (pushing existing rows 458-560 down to
# 459-561), then populate the new row with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 458..560 down by one, creating a blank row 458.
$ws.Rows.Item(458).Insert()

# Fill in the new row 458 with the new record.
$ws.Range("A458").Value = 4
$ws.Range("B458").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C458").Value = "Los Lagos"
$ws.Range("D458").Value = 45244
$ws.Range("E458").Value = 10
$ws.Range("F458").Value = 100114014
$ws.Range("G458").Value = "Betarraga"
$ws.Range("H458").Value = "Sin especificar"
$ws.Range("I458").Value = "Primera"
$ws.Range("J458").Value = 1200
$ws.Range("K458").Value = 1000
$ws.Range("L458").Value = 1100
$ws.Range("M458").Value = 1050
$ws.Range("N458").Value = "$/paquete 5 unidades"
$ws.Range("O458").Value = "Región Metropolitana"
$ws.Range("P458").Value = 210
$ws.Range("Q458").Value = 5
$ws.Range("R458").Value = "Hortaliza"
